# Add per-date sequence tags ("#N") to the experiment-name strings in
# column A. The tag reflects the Nth experiment recorded for that date,
# e.g. "2023-11-09 wastewater 0.485 gL - reflux valve open" (2nd entry
# for 2023-11-09) becomes "2023-11-09 #2 wastewater 0.485 gL - reflux
# valve open".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "2023-11-08 #1 tap water - reflux valve open"
    3  = "2023-11-09 #1 tap water - reflux valve open"
    4  = "2023-11-09 #2 wastewater 0.485 gL - reflux valve open"
    5  = "2023-11-09 #3 wastewater 1.297 gL - reflux valve open"
    6  = "2023-11-09 #4 wastewater 1.708 gL - reflux valve open"
    7  = "2024-07-11 #1 tap water - reflux valve closed"
    8  = "2024-07-12 #1 tap water - reflux valve open"
    9  = "2024-07-15 #1 wastewater 0.500 gL - reflux valve open"
    10 = "2024-07-16 #1 wastewater 0.250 gL - reflux valve open"
    11 = "2024-07-16 #2 wastewater 0.125 gL - reflux valve open"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}

# Reflect the final active-cell position recorded in the workbook after
# the edit.
$ws.Range("A12").Select()
